$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell {
    param($row, $col, $value)
    $cell = $ws.Cells.Item($row, $col)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = $origStyle
}

Set-TextCell 2 4 "70.157.83"
Set-TextCell 2 5 "  -1.34%  "

Set-TextCell 3 4 "3.738.72"
Set-TextCell 3 5 "  +0.47%  "

Set-TextCell 4 5 "  +0.05%  "

Set-TextCell 5 4 "623.44"
Set-TextCell 5 5 "  +0.50%  "

Set-TextCell 6 4 "180.58"
Set-TextCell 6 5 "  -1.36%  "

Set-TextCell 7 4 "3.735.61"
Set-TextCell 7 5 "  +0.45%  "

Set-TextCell 8 5 "  +0.06%  "

Set-TextCell 9 5 "  -1.45%  "

Set-TextCell 10 5 "  +2.19%  "

Set-TextCell 11 4 "6.30"
Set-TextCell 11 5 "  -5.00%  "

Set-TextCell 12 4 "0.488"
Set-TextCell 12 5 "  -3.08%  "

Set-TextCell 13 4 "40.94"
Set-TextCell 13 5 "  +0.53%  "

Set-TextCell 14 5 "  +1.57%  "

Set-TextCell 15 4 "4.361.79"
Set-TextCell 15 5 "  +0.48%  "

Set-TextCell 16 4 "3.736.04"
Set-TextCell 16 5 "  +0.61%  "

Set-TextCell 17 4 "70.164.96"
Set-TextCell 17 5 "  -1.35%  "

Set-TextCell 18 5 "  -1.25%  "

Set-TextCell 19 4 "7.61"
Set-TextCell 19 5 "  +0.93%  "

Set-TextCell 20 4 "16.85"
Set-TextCell 20 5 "  -0.52%  "

Set-TextCell 21 4 "506.32"
Set-TextCell 21 5 "  -2.57%  "

Set-TextCell 22 4 "9.35"
Set-TextCell 22 5 "  +0.15%  "

Set-TextCell 23 4 "0.727"
Set-TextCell 23 5 "  -2.56%  "

Set-TextCell 24 4 "2.59"
Set-TextCell 24 5 "  +1.48%  "

Set-TextCell 25 4 "86.72"
Set-TextCell 25 5 "  -2.44%  "

Set-TextCell 26 2 "RenderToken"
Set-TextCell 26 3 "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextCell 26 4 "11.46"
Set-TextCell 26 5 "  +1.97%  "

Set-TextCell 27 2 "InternetComputer(DFINITY)"
Set-TextCell 27 3 "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
Set-TextCell 27 4 "13.18"
Set-TextCell 27 5 "  -3.01%  "

Set-TextCell 28 5 "  +22.93%  "

Set-TextCell 29 5 "  -0.15%  "

Set-TextCell 30 5 "  -2.41%  "

Set-TextCell 31 4 "2.95"
Set-TextCell 31 5 "  +1.58%  "

Set-TextCell 32 4 "7.98"
Set-TextCell 32 5 "  -2.93%  "

Set-TextCell 33 4 "31.23"
Set-TextCell 33 5 "  -2.26%  "

Set-TextCell 34 4 "0.116"
Set-TextCell 34 5 "  -0.56%  "

Set-TextCell 35 4 "1.00"
Set-TextCell 35 5 "  +0.07%  "

Set-TextCell 36 4 "1.07"
Set-TextCell 36 5 "  +2.34%  "

Set-TextCell 37 5 "  +0.80%  "

Set-TextCell 38 5 "  +1.94%  "

Set-TextCell 39 4 "0.338"
Set-TextCell 39 5 "  -2.54%  "

Set-TextCell 40 4 "2.11"
Set-TextCell 40 5 "  -6.87%  "

Set-TextCell 41 4 "50.44"
Set-TextCell 41 5 "  -2.39%  "

Set-TextCell 42 4 "45.78"
Set-TextCell 42 5 "  +1.55%  "

Set-TextCell 43 4 "429.10"
Set-TextCell 43 5 "  -1.82%  "

Set-TextCell 44 5 "  -1.36%  "

Set-TextCell 45 4 "2.87"
Set-TextCell 45 5 "  -0.30%  "

Set-TextCell 46 4 "3.008.29"
Set-TextCell 46 5 "  -4.58%  "

Set-TextCell 47 5 "  -1.24%  "

Set-TextCell 48 4 "27.51"
Set-TextCell 48 5 "  -2.91%  "

Set-TextCell 50 2 "Monero"
Set-TextCell 50 3 "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextCell 50 4 "137.07"
Set-TextCell 50 5 "  -2.58%  "

Set-TextCell 51 2 "ThetaToken"
Set-TextCell 51 3 "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
Set-TextCell 51 4 "2.53"
Set-TextCell 51 5 "  +1.90%  "
